$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 250003470
$ws.Range("I64").Value = 4632.6665
$ws.Range("J64").Value = 1000000000
$ws.Range("K64").Value = 4632.6665
$ws.Range("L64").Value = 1000000000
$ws.Range("M64").Value = -4384.6665
$ws.Range("N64").Value = -1000000496
$ws.Range("H67").Value = 250003470
$ws.Range("I67").Value = 4632.6665
$ws.Range("J67").Value = 1000000000
$ws.Range("K67").Value = 4632.6665
$ws.Range("L67").Value = 1000000000
$ws.Range("M67").Value = -3774.6665
$ws.Range("N67").Value = -1000001716
$ws.Range("H88").Value = 8337042
$ws.Range("I88").Value = 12502625
$ws.Range("K88").Value = 12502625
$ws.Range("M88").Value = -12502219
$ws.Range("H91").Value = 8337042
$ws.Range("I91").Value = 12502625
$ws.Range("K91").Value = 12502625
$ws.Range("M91").Value = -12501221
$ws.Range("H115").Value = 320.2857
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H127").Value = 331.125
$ws.Range("I127").Value = 330.42856
$ws.Range("K127").Value = 991.28568
$ws.Range("M127").Value = 3968.71432
$ws.Range("H138").Value = 2771.7407
$ws.Range("J138").Value = 2940.1633
$ws.Range("L138").Value = 8820.4899
$ws.Range("N138").Value = -19100.4899
$ws.Range("H141").Value = 2180.08
$ws.Range("I141").Value = 2090.1
$ws.Range("K141").Value = 6270.299999999999
$ws.Range("M141").Value = -1090.299999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 200
$ws.Range("I10").Value = 200
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 200
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -30
$ws.Range("N10").ClearContents()
$ws.Range("H12").Value = 1990
$ws.Range("I12").Value = 1990
$ws.Range("K12").Value = 1990
$ws.Range("M12").Value = -1817
$ws.Range("H32").Value = 3479.3704
$ws.Range("I32").Value = 3479.3704
$ws.Range("K32").Value = 3479.3704
$ws.Range("M32").Value = -3192.3704
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 117653860
$ws.Range("I94").Value = 181827780
$ws.Range("J94").Value = 1664.5
$ws.Range("K94").Value = 181827780
$ws.Range("L94").Value = 1664.5
$ws.Range("M94").Value = -181827329
$ws.Range("N94").Value = -2566.5
$ws.Range("H128").Value = 11457.6
$ws.Range("I128").Value = 11457.6
$ws.Range("K128").Value = 34372.8
$ws.Range("M128").Value = -31882.8
$ws.Range("H134").Value = 2173.5789
$ws.Range("I134").Value = 1520.4348
$ws.Range("J134").Value = 3175.0667
$ws.Range("K134").Value = 4561.3044
$ws.Range("L134").Value = 9525.2001
$ws.Range("M134").Value = -2026.3044
$ws.Range("N134").Value = -14595.2001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5211.5854
$ws.Range("I31").Value = 3406.36
$ws.Range("K31").Value = 3406.36
$ws.Range("M31").Value = -3111.36
$ws.Range("H34").Value = 5211.5854
$ws.Range("I34").Value = 3406.36
$ws.Range("K34").Value = 3406.36
$ws.Range("M34").Value = -3204.36
$ws.Range("H58").Value = 1759
$ws.Range("I58").Value = 842.46155
$ws.Range("K58").Value = 842.46155
$ws.Range("M58").Value = -639.46155
$ws.Range("H93").Value = 6942.4
$ws.Range("I93").Value = 6942.4
$ws.Range("K93").Value = 6942.4
$ws.Range("M93").Value = -5070.4
$ws.Range("H132").Value = 11499454
$ws.Range("I132").Value = 3225.6667
$ws.Range("J132").Value = 30311464
$ws.Range("K132").Value = 9677.000100000001
$ws.Range("L132").Value = 90934392
$ws.Range("M132").Value = -7147.000100000001
$ws.Range("N132").Value = -90939452
$ws.Range("H134").Value = 4109.625
$ws.Range("I134").Value = 3383.6
$ws.Range("K134").Value = 10150.8
$ws.Range("M134").Value = -7615.799999999999
$ws.Range("H136").Value = 1759
$ws.Range("I136").Value = 842.46155
$ws.Range("K136").Value = 2527.38465
$ws.Range("M136").Value = 22.61535000000003
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 852.44446
$ws.Range("J5").Value = 872.5
$ws.Range("L5").Value = 2617.5
$ws.Range("N5").Value = -2841.5
$ws.Range("H68").Value = 16681908
$ws.Range("I68").Value = 25000
$ws.Range("J68").Value = 20013288
$ws.Range("K68").Value = 75000
$ws.Range("L68").Value = 60039864
$ws.Range("M68").Value = -74189
$ws.Range("N68").Value = -60041486
$ws.Range("H71").Value = 16681908
$ws.Range("I71").Value = 25000
$ws.Range("J71").Value = 20013288
$ws.Range("K71").Value = 225000
$ws.Range("L71").Value = 180119592
$ws.Range("M71").Value = -220944
$ws.Range("N71").Value = -180127704
$ws.Range("H113").Value = 1738.4
$ws.Range("I113").Value = 899
$ws.Range("J113").Value = 1948.25
$ws.Range("K113").Value = 2697
$ws.Range("L113").Value = 5844.75
$ws.Range("M113").Value = -527
$ws.Range("N113").Value = -10184.75
$ws.Range("H131").Value = 11365822
$ws.Range("J131").Value = 2261
$ws.Range("L131").Value = 6783
$ws.Range("N131").Value = -16863
$ws.Range("H132").Value = 2894.7
$ws.Range("J132").Value = 3714.1428
$ws.Range("L132").Value = 33427.2852
$ws.Range("N132").Value = -38487.2852
$ws.Range("H135").Value = 852.44446
$ws.Range("J135").Value = 872.5
$ws.Range("L135").Value = 7852.5
$ws.Range("N135").Value = -12922.5
$ws.Range("H139").Value = 5863.5557
$ws.Range("I139").Value = 3989.625
$ws.Range("K139").Value = 11968.875
$ws.Range("M139").Value = -6828.875
$ws.Range("H140").Value = 6540
$ws.Range("I140").Value = 2464.4736
$ws.Range("K140").Value = 7393.4208
$ws.Range("M140").Value = -2213.4208
$ws.Range("H141").Value = 8394
$ws.Range("I141").Value = 8394
$ws.Range("K141").Value = 25182
$ws.Range("M141").Value = -20002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 10.5
$ws.Range("I2").Value = 9
$ws.Range("K2").Value = 9
$ws.Range("M2").Value = 104
$ws.Range("H38").Value = 671675
$ws.Range("J38").Value = 671675
$ws.Range("L38").Value = 671675
$ws.Range("N38").Value = -672601
$ws.Range("H122").Value = 4052884.2
$ws.Range("I122").Value = 6996911
$ws.Range("K122").Value = 20990733
$ws.Range("M122").Value = -20988283
$ws.Range("H132").Value = 2489.0908
$ws.Range("I132").Value = 2359.125
$ws.Range("J132").Value = 2835.6667
$ws.Range("K132").Value = 7077.375
$ws.Range("L132").Value = 8507.000100000001
$ws.Range("M132").Value = -4547.375
$ws.Range("N132").Value = -13567.0001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 4871.75
$ws.Range("I31").Value = 395
$ws.Range("J31").Value = 12333
$ws.Range("K31").Value = 395
$ws.Range("L31").Value = 12333
$ws.Range("M31").Value = -147
$ws.Range("N31").Value = -12829
$ws.Range("H40").Value = 1499
$ws.Range("I40").Value = 1499
$ws.Range("K40").Value = 1499
$ws.Range("M40").Value = -1363
$ws.Range("H104").Value = 39666.332
$ws.Range("J104").Value = 39666.332
$ws.Range("L104").Value = 39666.332
$ws.Range("N104").Value = -46654.332
$ws.Range("H133").Value = 78968
$ws.Range("J133").Value = 78968
$ws.Range("L133").Value = 78968
$ws.Range("N133").Value = -84028
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5540.706
$ws.Range("I81").Value = 5219.3
$ws.Range("J81").Value = 5999.857
$ws.Range("K81").Value = 10438.6
$ws.Range("L81").Value = 11999.714
$ws.Range("M81").Value = -9377.6
$ws.Range("N81").Value = -14121.714
$ws.Range("H84").Value = 5540.706
$ws.Range("I84").Value = 5219.3
$ws.Range("J84").Value = 5999.857
$ws.Range("K84").Value = 52193
$ws.Range("L84").Value = 59998.57
$ws.Range("M84").Value = -46889
$ws.Range("N84").Value = -70606.57000000001
$ws.Range("H126").Value = 2260.125
$ws.Range("I126").Value = 2297.8572
$ws.Range("K126").Value = 6893.571599999999
$ws.Range("M126").Value = -4423.571599999999
